# Updates the cryptos list (Price / Volume(1h) columns, plus a few
# re-ordered coin rows) to match the latest GitHub Actions scrape.
# Note: some Price values look like plain numbers (e.g. "380.58"), so a
# leading apostrophe is used to force Excel to keep them as literal text
# instead of auto-converting them to numeric values (matching the
# original inlineStr/text storage of the Price column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.047.60'
$ws.Range('E2').Value = '  -0.27%  '

$ws.Range('D3').Value = '2.956.03'
$ws.Range('E3').Value = '  +0.42%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '''380.58'
$ws.Range('E5').Value = '  +1.03%  '

$ws.Range('D6').Value = '''102.09'
$ws.Range('E6').Value = '  -0.70%  '

$ws.Range('E7').Value = '  +1.76%  '

$ws.Range('E8').Value = '  +0.03%  '

$ws.Range('D9').Value = '''0.587'
$ws.Range('E9').Value = '  +0.53%  '

$ws.Range('D10').Value = '''36.50'
$ws.Range('E10').Value = '  -0.89%  '

$ws.Range('E11').Value = '  -0.84%  '

$ws.Range('D12').Value = '''0.0851'
$ws.Range('E12').Value = '  +1.63%  '

$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '3.421.30'
$ws.Range('E13').Value = '  +0.45%  '

$ws.Range('D14').Value = '''18.37'
$ws.Range('E14').Value = '  +2.31%  '

$ws.Range('B15').Value = 'Uniswap'
$ws.Range('C15').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D15').Value = '''12.34'
$ws.Range('E15').Value = '  +73.53%  '

$ws.Range('D16').Value = '''7.74'
$ws.Range('E16').Value = '  +5.20%  '

$ws.Range('D17').Value = '2.964.04'
$ws.Range('E17').Value = '  +1.56%  '

$ws.Range('D18').Value = '''1.00'
$ws.Range('E18').Value = '  +3.36%  '

$ws.Range('D19').Value = '51.109.35'
$ws.Range('E19').Value = '  +0.01%  '

$ws.Range('D20').Value = '''3.09'
$ws.Range('E20').Value = '  -2.44%  '

$ws.Range('D21').Value = '''12.36'
$ws.Range('E21').Value = '  -1.92%  '

$ws.Range('D22').Value = '0.0₃0959'
$ws.Range('E22').Value = '  +0.35%  '

$ws.Range('D23').Value = '''3.34'
$ws.Range('E23').Value = '  +16.12%  '

$ws.Range('D24').Value = '''269.06'
$ws.Range('E24').Value = '  +2.19%  '

$ws.Range('D25').Value = '''69.72'
$ws.Range('E25').Value = '  +2.17%  '

$ws.Range('E26').Value = '  -2.25%  '

$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('E28').Value = '  -1.04%  '

$ws.Range('D29').Value = '''25.89'
$ws.Range('E29').Value = '  +0.82%  '

$ws.Range('D30').Value = '''7.04'
$ws.Range('E30').Value = '  -11.36%  '

$ws.Range('E31').Value = '  -3.71%  '

$ws.Range('D32').Value = '''10.41'
$ws.Range('E32').Value = '  +5.55%  '

$ws.Range('E33').Value = '  +5.41%  '

$ws.Range('B34').Value = 'InjectiveProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D34').Value = '''34.36'
$ws.Range('E34').Value = '  +0.34%  '

$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''51.15'
$ws.Range('E35').Value = '  +0.60%  '

$ws.Range('D36').Value = '''0.0436'
$ws.Range('E36').Value = '  -4.75%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('E38').Value = '  +8.90%  '

$ws.Range('E39').Value = '  +1.93%  '

$ws.Range('E40').Value = '  +1.33%  '

$ws.Range('E42').Value = '  -3.18%  '

$ws.Range('D43').Value = '''124.59'
$ws.Range('E43').Value = '  +2.33%  '

$ws.Range('E44').Value = '  +10.01%  '

$ws.Range('D45').Value = '''21.60'
$ws.Range('E45').Value = '  +2.55%  '

$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').Value = '''0.274'
$ws.Range('E46').Value = '  -0.02%  '

$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.066.36'
$ws.Range('E47').Value = '  +3.34%  '

$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '''2.01'
$ws.Range('E48').Value = '  -1.88%  '

$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '''2.36'
$ws.Range('E49').Value = '  +1.18%  '

$ws.Range('E50').Value = '  -8.27%  '

$ws.Range('E51').Value = '  +6.40%  '
